$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17, pushing the existing rows 17-50 down to 18-51
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with the new weekly record
$ws.Range("A17").Value = 10
$ws.Range("B17").Value = "Vega Modelo de Temuco"
$ws.Range("C17").Value = "La Araucanía"
$ws.Range("D17").Value = 44525
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = 300000001
$ws.Range("G17").Value = "Rabanito"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 90
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 7000
$ws.Range("M17").Value = 6444
$ws.Range("N17").Value = "$/docena de paquetes"
$ws.Range("O17").Value = "Provincia de Cautín"
$ws.Range("P17").Value = 537
$ws.Range("Q17").Value = 12
$ws.Range("R17").Value = "Hortaliza"
